$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: extend header sequence with two new values (columns P and Q = 14, 15)
$ws.Cells.Item(1, 16).Value = 14   # P1
$ws.Cells.Item(1, 17).Value = 15   # Q1

# Copy the header style (bold, centered, bordered -> style "1") from O1 onto the
# two newly added header cells so they match the rest of row 1.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# Rows 2-25: swap values in columns I/K and M/O, and add new columns P, Q = 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2 (was 1)
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1 (was 2)
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2 (was 1)
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1 (was 2)
    $ws.Cells.Item($r, 16).Value = 2   # P (new)
    $ws.Cells.Item($r, 17).Value = 2   # Q (new)
}
